$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: C5 in_freq -> out of range divisor value, changes from 55 to 25000 ---
$ws.Range("C5").Value = 25000

# --- Row 6: CNT_VAL formulas, updated algorithm (shared formula D6:G6) ---
$ws.Range("C6").Formula = '=IF(AND(((C4/(2*$C$5))-1)<(2^$E$5),((C4/(2*$C$5))-1)>=0),ROUND(((C4/(2*$C$5))-1),0),"Impossible")'
$ws.Range("D6:G6").Formula = '=IF(AND(((D4/(2*$C$5))-1)<(2^$E$5),((D4/(2*$C$5))-1)>=0),ROUND(((D4/(2*$C$5))-1),0),"Impossible")'

# --- Row 7: real_out_freq formulas, updated to use CNT_VAL ---
$ws.Range("C7").Formula = '=C4/((C6+1)*2)'
$ws.Range("D7:G7").Formula = '=D4/((D6+1)*2)'

# --- Row 8: {USESCLK,CLKSRC} binary values - left-rotate C..G ---
$ws.Range("C8").Formula = '=DEC2BIN(0,3)'
$ws.Range("D8").Formula = '=DEC2BIN(1,3)'
$ws.Range("E8").Formula = '=DEC2BIN(2,3)'
$ws.Range("F8").Formula = '=DEC2BIN(3,3)'
$ws.Range("G8").Formula = '=DEC2BIN(4,3)'

# --- New rows 13/14: resolution/out_freq example computed with fixed CNT_VAL ---
$ws.Range("B13").Value = "CNT_VAL"
$ws.Range("C13").Value = 2
$ws.Range("B14").Value = "out_freq"
$ws.Range("C14").Formula = '=C4/(($C$13+1)*2)'
$ws.Range("D14:G14").Formula = '=D4/(($C$13+1)*2)'

# --- Selection moves to I7 ---
$ws.Range("I7").Select() | Out-Null
